$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 339.15384
$ws.Range("I2").Value = 257.375
$ws.Range("J2").Value = 470
$ws.Range("K2").Value = 257.375
$ws.Range("L2").Value = 470
$ws.Range("M2").Value = -144.375
$ws.Range("N2").Value = -696
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 130.92308
$ws.Range("I33").Value = 133.58333
$ws.Range("K33").Value = 133.58333
$ws.Range("M33").Value = 95.41667000000001
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1930.5
$ws.Range("I40").Value = 1862.875
$ws.Range("K40").Value = 1862.875
$ws.Range("M40").Value = -1687.875
# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 166
$ws.Range("I55").Value = 166
$ws.Range("K55").Value = 166
$ws.Range("M55").Value = 48
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1488.0217
$ws.Range("I98").Value = 1206.3889
$ws.Range("K98").Value = 1206.3889
$ws.Range("M98").Value = 291.6111000000001
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1488.0217
$ws.Range("I122").Value = 1206.3889
$ws.Range("K122").Value = 3619.1667
$ws.Range("M122").Value = -1169.1667
# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 450.5
$ws.Range("I125").Value = 450.5
$ws.Range("K125").Value = 4054.5
$ws.Range("M125").Value = -1594.5
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1080.1333
$ws.Range("I129").Value = 393.5
$ws.Range("J129").Value = 1112.0698
$ws.Range("K129").Value = 1180.5
$ws.Range("L129").Value = 3336.2094
$ws.Range("M129").Value = 3819.5
$ws.Range("N129").Value = -13336.2094
# Row 139 (Leve Item ID 42306)
$ws.Range("H139").Value = 67214.78
$ws.Range("J139").Value = 67214.78
$ws.Range("L139").Value = 67214.78
$ws.Range("N139").Value = -77494.78
# Row 140 (Leve Item ID 42459)
$ws.Range("H140").Value = 83900
$ws.Range("J140").Value = 83900
$ws.Range("L140").Value = 83900
$ws.Range("N140").Value = -94260

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1531.9333
$ws.Range("J45").Value = 1689
$ws.Range("L45").Value = 1689
$ws.Range("N45").Value = -2443
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 951.8788
$ws.Range("I74").Value = 557.5
$ws.Range("K74").Value = 557.5
$ws.Range("M74").Value = 316.5
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 951.8788
$ws.Range("I77").Value = 557.5
$ws.Range("K77").Value = 2787.5
$ws.Range("M77").Value = 1580.5
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 1239.3889
$ws.Range("I102").Value = 1221.1177
$ws.Range("K102").Value = 1221.1177
$ws.Range("M102").Value = 400.8823
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1227.1724
$ws.Range("I122").Value = 939.5833
$ws.Range("K122").Value = 2818.7499
$ws.Range("M122").Value = -368.7498999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1540.1428
$ws.Range("I99").Value = 1471
$ws.Range("K99").Value = 1471
$ws.Range("M99").Value = 27
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 1986.4193
$ws.Range("I105").Value = 1983.0358
$ws.Range("J105").Value = 2018
$ws.Range("K105").Value = 1983.0358
$ws.Range("L105").Value = 2018
$ws.Range("M105").Value = -236.0358000000001
$ws.Range("N105").Value = -5512

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 3 (Leve Item ID 3763)
$ws.Range("H3").Value = 6740
$ws.Range("J3").Value = 6740
$ws.Range("L3").Value = 6740
$ws.Range("N3").Value = -6966
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1515.7693
$ws.Range("I107").Value = 1246
$ws.Range("J107").Value = 2999.5
$ws.Range("K107").Value = 1246
$ws.Range("L107").Value = 2999.5
$ws.Range("M107").Value = 674
$ws.Range("N107").Value = -6839.5
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 3351.4
$ws.Range("I122").Value = 1915
$ws.Range("K122").Value = 5745
$ws.Range("M122").Value = -3295

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 56 (Leve Item ID 10146)
$ws.Range("H56").Value = 6730.6313
$ws.Range("I56").Value = 6730.6313
$ws.Range("K56").Value = 6730.6313
$ws.Range("M56").Value = -6200.6313
# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 496
$ws.Range("J86").Value = 145
$ws.Range("L86").Value = 435
$ws.Range("N86").Value = -2807
# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 496
$ws.Range("J89").Value = 145
$ws.Range("L89").Value = 1305
$ws.Range("N89").Value = -13161

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1067.2916
$ws.Range("I122").Value = 1019.6923
$ws.Range("K122").Value = 3059.0769
$ws.Range("M122").Value = -609.0769
# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 33711.555
$ws.Range("J123").Value = 33711.555
$ws.Range("L123").Value = 33711.555
$ws.Range("N123").Value = -38611.555
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3773826.2
$ws.Range("I126").Value = 4633373.5
$ws.Range("J126").Value = 335638
$ws.Range("K126").Value = 13900120.5
$ws.Range("L126").Value = 1006914
$ws.Range("M126").Value = -13897650.5
$ws.Range("N126").Value = -1011854

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3557.625
$ws.Range("I7").Value = 3326.8333
$ws.Range("K7").Value = 3326.8333
$ws.Range("M7").Value = -3214.8333
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 3446.8572
$ws.Range("I16").Value = 4130.8823
$ws.Range("J16").Value = 539.75
$ws.Range("K16").Value = 4130.8823
$ws.Range("L16").Value = 539.75
$ws.Range("M16").Value = -3960.8823
$ws.Range("N16").Value = -879.75
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 450.58334
$ws.Range("J55").Value = 494.91666
$ws.Range("L55").Value = 494.91666
$ws.Range("N55").Value = -840.91666
# Row 63 (Leve Item ID 12006)
$ws.Range("H63").Value = 46985
$ws.Range("J63").Value = 46985
$ws.Range("L63").Value = 46985
$ws.Range("N63").Value = -48483
# Row 66 (Leve Item ID 12006)
$ws.Range("H66").Value = 46985
$ws.Range("J66").Value = 46985
$ws.Range("L66").Value = 140955
$ws.Range("N66").Value = -148443
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2679.2354
$ws.Range("I122").Value = 2634.3845
$ws.Range("K122").Value = 7903.1535
$ws.Range("M122").Value = -5453.1535
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3557.625
$ws.Range("I126").Value = 3326.8333
$ws.Range("K126").Value = 9980.499899999999
$ws.Range("M126").Value = -7510.499899999999
# Row 131 (Leve Item ID 35466)
$ws.Range("H131").Value = 70000
$ws.Range("J131").Value = 70000
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1118
$ws.Range("I100").Value = 1022.5
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 2045
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1504
$ws.Range("N100").Value = -4082
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 808.2727
$ws.Range("I107").Value = 699.7778
$ws.Range("K107").Value = 2099.3334
$ws.Range("M107").Value = -179.3334
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 54920.867
$ws.Range("I122").Value = 58664.855
$ws.Range("K122").Value = 175994.565
$ws.Range("M122").Value = -173544.565
